# week 10 sum 2022 updates
# Adds a new "Week 28" column (AC) with this week's innings-count totals
# for the players who pitched, mirroring the existing Week N columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("AC1").Value = "Week 28"

# New week's totals (only players who played that week get a value)
$ws.Range("AC4").Value = 7
$ws.Range("AC6").Value = 7
$ws.Range("AC7").Value = 5
$ws.Range("AC8").Value = 3.5
$ws.Range("AC9").Value = 1.5

# Leave the selection where the user ended up after entering the data
[void]$ws.Range("AD12").Select()
